# Apply updated "Median (5%, 95% quantiles)" values in column B.
# The lower bound of each range is replaced with the true 5% quantile
# value (previously it had erroneously duplicated the median value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dash = [char]0x2013

$updates = @{
    2  = "0.05 (0" + $dash + "0.23)"
    3  = "0.06 (0" + $dash + "0.28)"
    4  = "0.06 (0" + $dash + "0.26)"
    9  = "62.84 (47.43" + $dash + "80.38)"
    11 = "113.08 (0.87" + $dash + "651.94)"
    12 = "340 (3" + $dash + "2983)"
    13 = "0.64 (0.33" + $dash + "0.87)"
    14 = "0.75 (0.29" + $dash + "0.9)"
    16 = "0.08 (0" + $dash + "0.52)"
    17 = "0.34 (0.07" + $dash + "0.67)"
    18 = "0.01 (0" + $dash + "0.13)"
    19 = "0.05 (0" + $dash + "0.21)"
    21 = "0.04 (0" + $dash + "0.24)"
    22 = "0.23 (0.02" + $dash + "0.6)"
    23 = "32.73 (2.1" + $dash + "84.77)"
    24 = "74.6 (5.39" + $dash + "296.57)"
    25 = "1024.54 (0.44" + $dash + "5263.56)"
    29 = "0.56 (0.08" + $dash + "0.9)"
    30 = "0.01 (0" + $dash + "0.55)"
    31 = "0.32 (0.04" + $dash + "0.8)"
    32 = "0.01 (0" + $dash + "0.07)"
    34 = "0.01 (0" + $dash + "0.04)"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
